$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8497.5
$ws.Range("I106").Value = 8497.5
$ws.Range("K106").Value = 8497.5
$ws.Range("M106").Value = -7866.5
$ws.Range("H137").Value = 2704499
$ws.Range("I137").Value = 3227361.8
$ws.Range("K137").Value = 9682085.399999999
$ws.Range("M137").Value = -9679535.399999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4386135
$ws.Range("I5").Value = 5263252
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 5263252
$ws.Range("L5").Value = 550
$ws.Range("M5").Value = -5263140
$ws.Range("N5").Value = -774
$ws.Range("H32").Value = 1808.76
$ws.Range("I32").Value = 1570.7191
$ws.Range("J32").Value = 3734.7273
$ws.Range("K32").Value = 1570.7191
$ws.Range("L32").Value = 3734.7273
$ws.Range("M32").Value = -1283.7191
$ws.Range("N32").Value = -4308.7273
$ws.Range("H74").Value = 4847681.5
$ws.Range("I74").Value = 5838332
$ws.Range("K74").Value = 5838332
$ws.Range("M74").Value = -5837458
$ws.Range("H77").Value = 4847681.5
$ws.Range("I77").Value = 5838332
$ws.Range("K77").Value = 29191660
$ws.Range("M77").Value = -29187292
$ws.Range("H132").Value = 79010.92999999999
$ws.Range("I132").Value = 50924.85
$ws.Range("J132").Value = 159256.86
$ws.Range("K132").Value = 152774.55
$ws.Range("L132").Value = 477770.58
$ws.Range("M132").Value = -150244.55
$ws.Range("N132").Value = -482830.58

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4386135
$ws.Range("I4").Value = 5263252
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 5263252
$ws.Range("L4").Value = 550
$ws.Range("M4").Value = -5263137
$ws.Range("N4").Value = -780
$ws.Range("H99").Value = 1020.04346
$ws.Range("I99").Value = 1055.2941
$ws.Range("J99").Value = 920.1667
$ws.Range("K99").Value = 1055.2941
$ws.Range("L99").Value = 920.1667
$ws.Range("M99").Value = 442.7058999999999
$ws.Range("N99").Value = -3916.1667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 76923580
$ws.Range("I22").Value = 142857580
$ws.Range("J22").Value = 591.6667
$ws.Range("K22").Value = 142857580
$ws.Range("L22").Value = 591.6667
$ws.Range("M22").Value = -142857230
$ws.Range("N22").Value = -1291.6667
$ws.Range("H31").Value = 1770.6571
$ws.Range("I31").Value = 1036.68
$ws.Range("J31").Value = 3605.6
$ws.Range("K31").Value = 1036.68
$ws.Range("L31").Value = 3605.6
$ws.Range("M31").Value = -741.6800000000001
$ws.Range("N31").Value = -4195.6
$ws.Range("H34").Value = 1770.6571
$ws.Range("I34").Value = 1036.68
$ws.Range("J34").Value = 3605.6
$ws.Range("K34").Value = 1036.68
$ws.Range("L34").Value = 3605.6
$ws.Range("M34").Value = -834.6800000000001
$ws.Range("N34").Value = -4009.6
$ws.Range("H58").Value = 27028648
$ws.Range("I58").Value = 34484404
$ws.Range("J58").Value = 1525.125
$ws.Range("K58").Value = 34484404
$ws.Range("L58").Value = 1525.125
$ws.Range("M58").Value = -34484201
$ws.Range("N58").Value = -1931.125
$ws.Range("H136").Value = 27028648
$ws.Range("I136").Value = 34484404
$ws.Range("J136").Value = 1525.125
$ws.Range("K136").Value = 103453212
$ws.Range("L136").Value = 4575.375
$ws.Range("M136").Value = -103450662
$ws.Range("N136").Value = -9675.375

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3566.389
$ws.Range("I118").Value = 378
$ws.Range("J118").Value = 4080.6453
$ws.Range("K118").Value = 1134
$ws.Range("L118").Value = 12241.9359
$ws.Range("M118").Value = 109
$ws.Range("N118").Value = -14727.9359
$ws.Range("H131").Value = 969.32184
$ws.Range("I131").Value = 485.7
$ws.Range("J131").Value = 1032.1299
$ws.Range("K131").Value = 1457.1
$ws.Range("L131").Value = 3096.3897
$ws.Range("M131").Value = 3582.9
$ws.Range("N131").Value = -13176.3897
$ws.Range("H132").Value = 927.3077
$ws.Range("I132").Value = 681.875
$ws.Range("J132").Value = 1320
$ws.Range("K132").Value = 6136.875
$ws.Range("L132").Value = 11880
$ws.Range("M132").Value = -3606.875
$ws.Range("N132").Value = -16940

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4142.1665
$ws.Range("J80").Value = 4175.091
$ws.Range("L80").Value = 4175.091
$ws.Range("N80").Value = -6171.091
$ws.Range("H83").Value = 4142.1665
$ws.Range("J83").Value = 4175.091
$ws.Range("L83").Value = 20875.455
$ws.Range("N83").Value = -30859.455
$ws.Range("H132").Value = 288542.56
$ws.Range("I132").Value = 500900
$ws.Range("J132").Value = 203599.6
$ws.Range("K132").Value = 1502700
$ws.Range("L132").Value = 610798.8
$ws.Range("M132").Value = -1500170
$ws.Range("N132").Value = -615858.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2700
$ws.Range("I40").Value = 2700
$ws.Range("K40").Value = 2700
$ws.Range("M40").Value = -2564
$ws.Range("H46").Value = 986.5714
$ws.Range("I46").Value = 976
$ws.Range("J46").Value = 1000.6667
$ws.Range("K46").Value = 976
$ws.Range("L46").Value = 1000.6667
$ws.Range("M46").Value = -788
$ws.Range("N46").Value = -1376.6667
$ws.Range("H100").Value = 51504
$ws.Range("I100").Value = 84508.336
$ws.Range("J100").Value = 1997.5
$ws.Range("K100").Value = 84508.336
$ws.Range("L100").Value = 1997.5
$ws.Range("M100").Value = -83967.336
$ws.Range("N100").Value = -3079.5
$ws.Range("H132").Value = 52902.465
$ws.Range("I132").Value = 32163.059
$ws.Range("J132").Value = 80023.234
$ws.Range("K132").Value = 96489.177
$ws.Range("L132").Value = 240069.702
$ws.Range("M132").Value = -93959.177
$ws.Range("N132").Value = -245129.702

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3165
$ws.Range("I62").Value = 2865.3333
$ws.Range("J62").Value = 3726.875
$ws.Range("K62").Value = 2865.3333
$ws.Range("L62").Value = 3726.875
$ws.Range("M62").Value = -2241.3333
$ws.Range("N62").Value = -4974.875
$ws.Range("H65").Value = 3165
$ws.Range("I65").Value = 2865.3333
$ws.Range("J65").Value = 3726.875
$ws.Range("K65").Value = 14326.6665
$ws.Range("L65").Value = 18634.375
$ws.Range("M65").Value = -11206.6665
$ws.Range("N65").Value = -24874.375
$ws.Range("H132").Value = 81825.39999999999
$ws.Range("I132").Value = 78202.84
$ws.Range("K132").Value = 234608.52
$ws.Range("M132").Value = -232078.52
